$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text (avoid Excel auto-converting
# numeric-looking strings like "0.9980" or "0.05300" into numbers and
# dropping significant trailing zeros / switching to exponent notation).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.909.89'
$ws.Range('E2').Value = '  +3.73%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.681.74'
$ws.Range('E3').Value = '  +2.58%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9965'
$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.98'
$ws.Range('E5').Value = '  +2.42%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9987'
$ws.Range('E6').Value = '  -0.16%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4636'
$ws.Range('E7').Value = '  -2.13%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2607'
$ws.Range('E8').Value = '  +1.73%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06158'
$ws.Range('E9').Value = '  +1.36%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.675.12'

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07007'
$ws.Range('E11').Value = '  -0.25%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.98'
$ws.Range('E12').Value = '  +3.62%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.374'
$ws.Range('E13').Value = '  +1.41%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5811'
$ws.Range('E14').Value = '  +0.62%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '75.65'
$ws.Range('E15').Value = '  +2.76%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.9983'
$ws.Range('E16').Value = '  -0.19%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9980'

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.870.49'
$ws.Range('E18').Value = '  +3.66%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000006733'
$ws.Range('E19').Value = '  +2.44%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.45'
$ws.Range('E20').Value = '  +2.06%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.889.13'
$ws.Range('E21').Value = '  +1.80%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.479'
$ws.Range('E22').Value = '  +3.94%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.688'
$ws.Range('E23').Value = '  +2.43%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.242'
$ws.Range('E24').Value = '  +0.45%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '134.15'
$ws.Range('E25').Value = '  +0.90%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.02'
$ws.Range('E26').Value = '  +1.01%  '

$ws.Range('E27').Value = '  +1.31%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.717'
$ws.Range('E28').Value = '  +5.22%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '104.89'
$ws.Range('E29').Value = '  +0.82%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.961'
$ws.Range('E30').Value = '  +1.65%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07713'
$ws.Range('E31').Value = '  +2.11%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.628'
$ws.Range('E32').Value = '  +2.84%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04374'
$ws.Range('E33').Value = '  +2.94%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.595'
$ws.Range('E34').Value = '  +0.87%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9578'
$ws.Range('E35').Value = '  +3.33%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6083'
$ws.Range('E36').Value = '  +2.95%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.9364'
$ws.Range('E37').Value = '  +5.42%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '109.54'
$ws.Range('E38').Value = '  +11.13%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.461'
$ws.Range('E39').Value = '  -4.20%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9972'
$ws.Range('E40').Value = '  -0.27%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.882'
$ws.Range('E41').Value = '  +6.89%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.01459'
$ws.Range('E42').Value = '  -2.09%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.068'
$ws.Range('E43').Value = '  +8.90%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.3736'
$ws.Range('E44').Value = '  +1.28%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1123'
$ws.Range('E45').Value = '  +2.23%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.196'
$ws.Range('E46').Value = '  +1.82%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05300'
$ws.Range('E47').Value = '  +2.03%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '30.98'
$ws.Range('E48').Value = '  +8.55%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.663'
$ws.Range('E49').Value = '  +7.67%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.215'
$ws.Range('E50').Value = '  +2.72%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9997'
$ws.Range('E51').Value = '  -0.16%  '
